# Petal_waterloss_time.xlsx edit
#
# Sheet1, rows 26-31: the "population" column (A) was left blank for the
# second batch of petal-loss readings even though every other row in the
# sheet is tagged "cojo" (the Cojo HQ population). Backfill column A for
# those six rows so they match the rest of the table, then leave the
# selection where the author ended up (B34) after doing the data entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 26..31) {
    $ws.Cells.Item($r, 1).Value = "cojo"
}

$ws.Range("B34").Select()
